$wb = $excel.ActiveWorkbook

# Add a new worksheet "Hydration_all" after the existing (last) sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Hydration_all"

# Header row.
$ws2.Range("A1").Value = "Size"
$ws2.Range("B1").Value = "HydrationEnergy_kJperMol"
$ws2.Range("C1").Value = "Hyd2"
$ws2.Range("D1").Value = "psize params run"

# Re-simulated data (psize params run): Size, raw HydrationEnergy_kJperMol,
# Hyd2 (previous HydrationEnergy values), and the new D column.
$rows = @(
    @(30,  -29853.523000000001, -25304.7703,          -28413.667369440001),
    @(40,  -32100,              -27099.15,             -77711.659220169997),
    @(50,  -43537.3,            -32838.676800000001,   -75392.540148309999),
    @(60,  -37961.97,           -37961.97,              -82534.183866730004),
    @(70,  -39400,              -39400,                 -83708.873504570001),
    @(80,  -19845.761299999998, -19845.761299999998,    12110.654760789999),
    @(90,  -20155.915000000001, -20155.915000000001,    $null),
    @(100, -33756.466999999997, -33756.466999999997,   -107927.2195149)
)

$arr = New-Object 'object[,]' 8,4
for ($i = 0; $i -lt 8; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}
$ws2.Range("A2:D9").Value = $arr

# Column D (the new simulation output) gets an explicit "General" number
# format - this is what produces the extra cellXfs entry in the saved file.
# Skip D8 (row 8 has no value in that column) so we don't materialize a
# spurious styled-but-empty cell there.
$ws2.Range("D1:D7").NumberFormat = "General"
$ws2.Range("D9").NumberFormat = "General"

# Match the source column width for the new data column.
$ws2.Columns("D").ColumnWidth = 9.21875

# The new sheet becomes the active / selected tab, with this cell selected.
$ws2.Range("Q13").Select() | Out-Null
